$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1708185053380783
$ws.Range("C2").Value = 0.6014234875444839
$ws.Range("J2").Value = 0.01779359430604982
$ws.Range("P2").Value = 0.1352313167259787
$ws.Range("S2").Value = 0.07473309608540925
$ws.Range("B3").Value = 0.005747126436781609
$ws.Range("C3").Value = 0.01149425287356322
$ws.Range("J3").Value = 0.02873563218390805
$ws.Range("P3").Value = 0.8045977011494253
$ws.Range("S3").Value = 0.1494252873563219
$ws.Range("J4").Value = 0.09433962264150944
$ws.Range("P4").Value = 0.7547169811320755
$ws.Range("S4").Value = 0.1509433962264151
$ws.Range("P5").Value = 0.75
$ws.Range("S5").Value = 0.25
$ws.Range("B6").Value = 0.0730593607305936
$ws.Range("D6").Value = 0.0228310502283105
$ws.Range("F6").Value = 0.0684931506849315
$ws.Range("J6").Value = 0.2328767123287671
$ws.Range("O6").Value = 0.0319634703196347
$ws.Range("Q6").Value = 0.2328767123287671
$ws.Range("R6").Value = 0.0410958904109589
$ws.Range("S6").Value = 0.2968036529680365
$ws.Range("B7").Value = 0.09239130434782608
$ws.Range("D7").Value = 0.02173913043478261
$ws.Range("F7").Value = 0.05434782608695652
$ws.Range("J7").Value = 0.1902173913043478
$ws.Range("O7").Value = 0.0108695652173913
$ws.Range("Q7").Value = 0.2010869565217391
$ws.Range("R7").Value = 0.07065217391304347
$ws.Range("S7").Value = 0.358695652173913
$ws.Range("B8").Value = 0.09677419354838709
$ws.Range("D8").Value = 0.02764976958525346
$ws.Range("E8").Value = 0.002304147465437788
$ws.Range("F8").Value = 0.06221198156682028
$ws.Range("J8").Value = 0.0944700460829493
$ws.Range("O8").Value = 0.03456221198156682
$ws.Range("Q8").Value = 0.2235023041474654
$ws.Range("R8").Value = 0.07142857142857142
$ws.Range("S8").Value = 0.3870967741935484
$ws.Range("B9").Value = 0.1118421052631579
$ws.Range("D9").Value = 0.02631578947368421
$ws.Range("F9").Value = 0.03947368421052631
$ws.Range("J9").Value = 0.125
$ws.Range("O9").Value = 0.03289473684210526
$ws.Range("Q9").Value = 0.1776315789473684
$ws.Range("R9").Value = 0.05263157894736842
$ws.Range("S9").Value = 0.4342105263157895
$ws.Range("B10").Value = 0.1066260472201066
$ws.Range("D10").Value = 0.02208682406702209
$ws.Range("E10").Value = 0.002284843869002285
$ws.Range("F10").Value = 0.07844630616907845
$ws.Range("J10").Value = 0.1142421934501142
$ws.Range("O10").Value = 0.01751713632901752
$ws.Range("Q10").Value = 0.246001523229246
$ws.Range("R10").Value = 0.06626047220106626
$ws.Range("S10").Value = 0.3465346534653465
$ws.Range("G11").Value = 0.1084745762711864
$ws.Range("J11").Value = 0.09152542372881356
$ws.Range("K11").Value = 0.1796610169491525
$ws.Range("L11").Value = 0.6
$ws.Range("S11").Value = 0.02033898305084746
$ws.Range("G12").Value = 0.7377049180327869
$ws.Range("J12").Value = 0.2185792349726776
$ws.Range("K12").Value = 0.00546448087431694
$ws.Range("L12").Value = 0.02185792349726776
$ws.Range("S12").Value = 0.01639344262295082
$ws.Range("G13").Value = 0.5833333333333334
$ws.Range("J13").Value = 0.3888888888888889
$ws.Range("S13").Value = 0.02777777777777778
$ws.Range("G14").Value = 1
$ws.Range("F15").Value = 0.01984126984126984
$ws.Range("H15").Value = 0.1785714285714286
$ws.Range("I15").Value = 0.05952380952380952
$ws.Range("J15").Value = 0.3095238095238095
$ws.Range("K15").Value = 0.08333333333333333
$ws.Range("M15").Value = 0.01587301587301587
$ws.Range("O15").Value = 0.03968253968253968
$ws.Range("S15").Value = 0.2936507936507937
$ws.Range("F16").Value = 0.02325581395348837
$ws.Range("H16").Value = 0.1534883720930233
$ws.Range("I16").Value = 0.05116279069767442
$ws.Range("J16").Value = 0.4465116279069767
$ws.Range("K16").Value = 0.1116279069767442
$ws.Range("M16").Value = 0.0186046511627907
$ws.Range("O16").Value = 0.06046511627906977
$ws.Range("S16").Value = 0.1348837209302326
$ws.Range("F17").Value = 0.0205607476635514
$ws.Range("H17").Value = 0.1775700934579439
$ws.Range("I17").Value = 0.06355140186915888
$ws.Range("J17").Value = 0.4485981308411215
$ws.Range("K17").Value = 0.102803738317757
$ws.Range("M17").Value = 0.0205607476635514
$ws.Range("O17").Value = 0.06728971962616823
$ws.Range("S17").Value = 0.09906542056074766
$ws.Range("F18").Value = 0.02684563758389262
$ws.Range("H18").Value = 0.1812080536912752
$ws.Range("I18").Value = 0.1073825503355705
$ws.Range("J18").Value = 0.4228187919463087
$ws.Range("K18").Value = 0.09395973154362416
$ws.Range("M18").Value = 0.006711409395973154
$ws.Range("O18").Value = 0.08053691275167785
$ws.Range("S18").Value = 0.08053691275167785
$ws.Range("F19").Value = 0.01371036846615253
$ws.Range("H19").Value = 0.1996572407883462
$ws.Range("I19").Value = 0.06598114824335904
$ws.Range("J19").Value = 0.4001713796058269
$ws.Range("K19").Value = 0.1071122536418166
$ws.Range("M19").Value = 0.01371036846615253
$ws.Range("N19").Value = 0.000856898029134533
$ws.Range("O19").Value = 0.08483290488431877
$ws.Range("S19").Value = 0.1139674378748929
